# Add a "Post" request row (with a wrapped JSON body column) to the
# REST API test-data sheet, mirroring the existing "Get" row's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting first (copy/paste of formats), values are assigned
#     afterwards so new shared-string entries land in the same order the
#     original author typed them in: Post, Body, then the JSON body. ---

# A3 / B3 get the existing bordered style plus vertical-centering.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").VerticalAlignment = -4108

$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").VerticalAlignment = -4108

# C3 keeps the plain bordered style used by C2.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# D3 gets the bordered style, vertically centered AND wrapped, to hold the
# multi-line JSON request body.
$ws.Range("A2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true

# New "Body" column header (D1), copying the look of the other headers.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# D2 stays empty but picks up the same bordered look as the rest of row 2.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# Now assign the cell values/text, in the order the shared strings were
# first introduced: "Post" (A3), "Body" (D1), the JSON body (D3).
$ws.Range("A3").Value = "Post"
$ws.Range("B3").Value = "https://reqres.in/"
$ws.Range("C3").Value = 201
$ws.Range("D1").Value = "Body"
$ws.Range("D3").Value = "{`n    ""name"": ""RestTest"",`n    ""job"": ""Org""`n}"

# Row 3 is taller to show the wrapped body text.
$ws.Rows.Item(3).RowHeight = 60

# Widen the new column so the body text is readable (target sheet stores
# width 24.140625; this host's column grid only lands on sixths of a
# character, so 23.33 is the input that rounds to the nearest XML width).
$ws.Columns.Item(4).ColumnWidth = 23.33

# Clear the clipboard marquee and move the active selection, matching the
# saved workbook state.
$excel.CutCopyMode = $false
$ws.Range("D9").Select()
